# adding a while loop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (model, invoice, qty, stock flag, cost, unit)
# derived from shifting/renumbering the original SO2311035 line items into
# the new SO2311043 invoice, then dropping the trailing two rows.
$models = @(
    "T006.407.11.033.00",
    "T006.407.11.053.00",
    "T006.407.16.033.00",
    "T094.210.11.111.00",
    "T099.407.11.048.00",
    "T116.617.11.057.01",
    "T120.417.11.041.01",
    "T120.417.11.091.01",
    "T126.010.11.013.00",
    "T126.010.22.013.01",
    "T126.010.36.013.00",
    "T126.207.11.013.00",
    "T137.407.11.041.00",
    "T137.407.11.051.00",
    "T137.410.11.041.00"
)

$qtys = @(20, 25, 25, 15, 5, 5, 5, 5, 5, 10, 3, 5, 5, 5, 2)
$costs = @(3127, 3127, 2832, 1652, 3953, 1888, 2743.5, 2743.5, 2006, 2212.5, 2006, 3481, 3274.5, 3274.5, 1740.5)

$invoice = "SO2311043"

$i = 0
$row = 2
while ($i -lt $models.Length) {
    $ws.Cells.Item($row, 1).Value = $models[$i]
    $ws.Cells.Item($row, 2).Value = $invoice
    $ws.Cells.Item($row, 3).Value = $qtys[$i]
    $ws.Cells.Item($row, 5).Value = $costs[$i]
    $i = $i + 1
    $row = $row + 1
}

# remove the two trailing rows that no longer exist on the new invoice
$ws.Rows("17:18").Delete()

# update the selection shown when the sheet is reopened
$ws.Range("F2:F16").Select()
